# Add two new "Questão" slides (51 and 52) at the end of the deck,
# cloned from the last existing "Questão" slide so they inherit the
# exact same layout/shapes/pictures, then overwrite their text content.

$p = $ppt.ActivePresentation

# The last slide in the deck is "Questão 50" - use it as the template
# for the two new slides so shape ids / xfrm / pictures all match.
$templateIndex = $p.Slides.Count
$template = $p.Slides.Item($templateIndex)

$slide51 = $template.Duplicate().Item(1)
$slide52 = $slide51.Duplicate().Item(1)

# ---------------------------------------------------------------
# Slide 51 - "Questão 51" (Transformação Digital / Dissertativa)
# ---------------------------------------------------------------
$slide51.Shapes.Item(2).TextFrame.TextRange.Text = "Questão 51"

$questao51 = @(
  "Introdução:",
  "A utilização de dados e tecnologias digitais para experimentação, validação de hipóteses e tomada de decisão é cada vez mais presente em diversas áreas de negócio. A Ciência de Dados e o uso da Transformação Digital tem permitido a melhoria das decisões tomadas com base em dados e evidências.",
  "",
  "Questão:",
  "Qual das opções abaixo melhor descreve uma das vantagens de se utilizar a Transformação Digital para tomar decisões baseadas em dados?",
  "",
  "A) Diminuição do custo das decisões;",
  "B) Amadurecimento da tecnologia usada;",
  "C) Maior velocidade no processamento dos dados;",
  "D) Maior precisão nas decisões. ",
  "",
  "Resposta Correta: D) Maior precisão nas decisões.",
  "Justificativa: A Transformação Digital tem permitido uma maior previsibilidade e precisão nos processos de tomada de decisão, pois permite racionalizar o processo de coleta, produção e análise de dados, em comparação aos mecanismos convencionais. A opção A) não está relacionada com a Transformação Digital, pois esta pode ajudar a diminuir custos, mas não é o seu principal objetivo. A opção B) não está relacionada com a Transformação Digital, pois esta não tem como objetivo amadurecer as tecnologias usadas. A opção C) está relacionada, pois uma das vantagens é exatamente a maior velocidade no processamento de dados, mas não é o principal objetivo. Já a opção D) descreve corretamente uma das vantagens da Transformação Digital para decisões baseadas em dados: maior precisão nas decisões."
) -join "`r"
$slide51.Shapes.Item(3).TextFrame.TextRange.Text = $questao51

$slide51.Shapes.Item(5).TextFrame.TextRange.Text = "Dissertativa"
$slide51.Shapes.Item(7).TextFrame.TextRange.Text = "Transformação Digital"
$slide51.Shapes.Item(9).TextFrame.TextRange.Text = "Complexo"

# ---------------------------------------------------------------
# Slide 52 - "Questão 52" (Transformação Digital / Completar as Lacunas)
# ---------------------------------------------------------------
$slide52.Shapes.Item(2).TextFrame.TextRange.Text = "Questão 52"

$questao52 = @(
  "Questão:",
  "",
  "A tecnologia blockchain é usada para melhorar a transparência e segurança em diversas iniciativas transformadoras. Blockchain fornece às empresas a capacidade de trabalhar com transações descentralizadas nos setores de supply chain, saúde e finanças, permitindo a execução de transações em uma rede de computadores relacionados, o que proporciona maior segurança e ______.",
  "",
  "(A) audibilidade",
  "(B) compartilhamento",
  "(C) escalabilidade",
  "(D) imutabilidade",
  "(E) visibilidade",
  "",
  "Resposta: D - imutabilidade. A blockchain permite que transações sejam executadas de forma mais segura e imutável, o que significa que os dados armazenados não podem ser alterados ou excluídos dos registros."
) -join "`r"
$slide52.Shapes.Item(3).TextFrame.TextRange.Text = $questao52

$slide52.Shapes.Item(5).TextFrame.TextRange.Text = "Completar as Lacunas"
$slide52.Shapes.Item(7).TextFrame.TextRange.Text = "Transformação Digital"
$slide52.Shapes.Item(9).TextFrame.TextRange.Text = "Fácil"

Write-Output "Slides now: $($p.Slides.Count)"
